$d = $word.ActiveDocument
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# 1) Update the text of the existing "For our project ..." paragraph first,
#    while it is still the 3rd paragraph (and before we insert new ones
#    above it, so the Find targets stay unambiguous).
# ---------------------------------------------------------------------------
$old1 = "If we were to create a regular employee class, it wouldn" + $apos + "t extend management employee. It would only extend the employee abstract class that gives it permission to create maintenance requests. Our facility secretary class"
$new1 = "If we were to create a lower level employee class, it wouldn" + $apos + "t extend management employee. However all employees have the capability to make maintenance requests as they should be through the Employee abstract class. Our facility secretary class"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "The facilities client in the view package is where methods would be called."
$new2 = "The facilities client in the view package is where methods are called and input is taken through the console. Our project has tests that test certain methods."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert 6 new paragraphs before that paragraph (which is still index 3).
# ---------------------------------------------------------------------------
$target = $d.Paragraphs.Item(3)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> new para 3 (will hold "Presentation ...")
$p = $d.Paragraphs.Item(3).Range
$p.Collapse(1)
$p.InsertAfter("Presentation is separate from business logic and lives in the View package.")
$p.Font.Bold = 1

$target = $d.Paragraphs.Item(4)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> blank para 4
$target = $d.Paragraphs.Item(5)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> new para 5 ("To run our application ...")
$p = $d.Paragraphs.Item(5).Range
$p.Collapse(1)
$p.InsertAfter("To run our application, click the " + $apos + "run as" + $apos + " button from the project menu and run as java application. Follow the instructions to add a facility or exit the program. Follow the instructions to check the usage rate of a facility also.")

$target = $d.Paragraphs.Item(6)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> blank para 6
$target = $d.Paragraphs.Item(7)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> new para 7 ("To run tests ...")
$p = $d.Paragraphs.Item(7).Range
$p.Collapse(1)
$p.InsertAfter("To run tests, go to FacilityTest.java and click run. All tests should pass based on the parameters you enter for the assert methods. You can test our usage of certain parameters as an example. All tests pass.")

$target = $d.Paragraphs.Item(8)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> blank para 8
$target = $d.Paragraphs.Item(9)
$target.Range.InsertParagraphBefore() | Out-Null   ; # -> new para 9 ("More Information on design decisions:")
$p = $d.Paragraphs.Item(9).Range
$p.Collapse(1)
$p.InsertAfter("More Informatio")
$p.Font.Bold = 1
$p.Collapse(0)
$p.InsertAfter("n on design decisions")
$p.Font.Bold = 1
$p.Collapse(0)
$p.InsertAfter(":")
$p.Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark from the end of the big paragraph to the end
#    of the first new ("Presentation ...") paragraph, matching the diff.
# ---------------------------------------------------------------------------
$presentationPara = $d.Paragraphs.Item(3)
$bmRange = $presentationPara.Range.Duplicate
$bmRange.Collapse(0)
$bmRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "Para $i [" $d.Paragraphs.Item($i).Range.Text "]"
}
